$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 240.77777
$ws.Range("I33").Value = 240.77777
$ws.Range("K33").Value = 240.77777
$ws.Range("M33").Value = -11.77777

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3862.7334
$ws.Range("I40").Value = 3494.1
$ws.Range("K40").Value = 3494.1
$ws.Range("M40").Value = -3319.1

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1828.5
$ws.Range("I43").Value = 1983.1666
$ws.Range("J43").Value = 1712.5
$ws.Range("K43").Value = 1983.1666
$ws.Range("L43").Value = 1712.5
$ws.Range("M43").Value = -1914.1666
$ws.Range("N43").Value = -1850.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3643.625
$ws.Range("I62").Value = 3153.6924
$ws.Range("K62").Value = 3153.6924
$ws.Range("M62").Value = -2529.6924

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3643.625
$ws.Range("I65").Value = 3153.6924
$ws.Range("K65").Value = 15768.462
$ws.Range("M65").Value = -12648.462

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2795.4614
$ws.Range("I106").Value = 2695.0833
$ws.Range("K106").Value = 2695.0833
$ws.Range("M106").Value = -2064.0833

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2783566
$ws.Range("I116").Value = 5896.5
$ws.Range("J116").Value = 8338905
$ws.Range("K116").Value = 5896.5
$ws.Range("L116").Value = 8338905
$ws.Range("M116").Value = -2454.5
$ws.Range("N116").Value = -8345789

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1611.7646
$ws.Range("I132").Value = 1492.5385
$ws.Range("K132").Value = 4477.6155
$ws.Range("M132").Value = -1947.6155

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 114217.664
$ws.Range("I61").Value = 3342.2
$ws.Range("K61").Value = 3342.2
$ws.Range("M61").Value = -3130.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 37500
$ws.Range("J64").Value = 55000
$ws.Range("L64").Value = 55000
$ws.Range("N64").Value = -55496

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H67").Value = 37500
$ws.Range("J67").Value = 55000
$ws.Range("L67").Value = 55000
$ws.Range("N67").Value = -56716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 114217.664
$ws.Range("I136").Value = 3342.2
$ws.Range("K136").Value = 10026.6
$ws.Range("M136").Value = -7476.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2151343.2
$ws.Range("I99").Value = 125791.125
$ws.Range("J99").Value = 4466259.5
$ws.Range("K99").Value = 125791.125
$ws.Range("L99").Value = 4466259.5
$ws.Range("M99").Value = -124293.125
$ws.Range("N99").Value = -4469255.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4168729
$ws.Range("J107").Value = 2544.1428
$ws.Range("L107").Value = 2544.1428
$ws.Range("N107").Value = -6384.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2636
$ws.Range("I31").Value = 2031.6296
$ws.Range("K31").Value = 2031.6296
$ws.Range("M31").Value = -1736.6296

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2636
$ws.Range("I34").Value = 2031.6296
$ws.Range("K34").Value = 2031.6296
$ws.Range("M34").Value = -1829.6296

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1591.6154
$ws.Range("I58").Value = 1676
$ws.Range("J58").Value = 1493.1666
$ws.Range("K58").Value = 1676
$ws.Range("L58").Value = 1493.1666
$ws.Range("M58").Value = -1473
$ws.Range("N58").Value = -1899.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4656.857
$ws.Range("I62").Value = 4804.8335
$ws.Range("K62").Value = 4804.8335
$ws.Range("M62").Value = -4180.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4656.857
$ws.Range("I65").Value = 4804.8335
$ws.Range("K65").Value = 24024.1675
$ws.Range("M65").Value = -20904.1675

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1411.6875
$ws.Range("I107").Value = 1088.8572
$ws.Range("J107").Value = 2028
$ws.Range("K107").Value = 1088.8572
$ws.Range("L107").Value = 2028
$ws.Range("M107").Value = 831.1428000000001
$ws.Range("N107").Value = -5868

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1591.6154
$ws.Range("I136").Value = 1676
$ws.Range("J136").Value = 1493.1666
$ws.Range("K136").Value = 5028
$ws.Range("L136").Value = 4479.4998
$ws.Range("M136").Value = -2478
$ws.Range("N136").Value = -9579.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 585.6667
$ws.Range("I92").Value = 588.5
$ws.Range("J92").Value = 580
$ws.Range("K92").Value = 1765.5
$ws.Range("L92").Value = 1740
$ws.Range("M92").Value = -517.5
$ws.Range("N92").Value = -4236

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1182.625
$ws.Range("I102").Value = 1173.0714
$ws.Range("K102").Value = 1173.0714
$ws.Range("M102").Value = 448.9286

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2872838
$ws.Range("I113").Value = 187370.17
$ws.Range("K113").Value = 187370.17
$ws.Range("M113").Value = -185200.17

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 16000
$ws.Range("J64").Value = 16000
$ws.Range("L64").Value = 16000
$ws.Range("N64").Value = -16450

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H67").Value = 16000
$ws.Range("J67").Value = 16000
$ws.Range("L67").Value = 16000
$ws.Range("N67").Value = -17560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 30000
$ws.Range("J21").Value = 30000
$ws.Range("L21").Value = 30000
$ws.Range("N21").Value = -30470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 259999.8
$ws.Range("J24").Value = 74999.75
$ws.Range("L24").Value = 74999.75
$ws.Range("N24").Value = -75459.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H35").Value = 30000
$ws.Range("J35").Value = 30000
$ws.Range("L35").Value = 30000
$ws.Range("N35").Value = -30580

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 3402062.2
$ws.Range("I100").Value = 3968940
$ws.Range("K100").Value = 7937880
$ws.Range("M100").Value = -7937339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 187751.33
$ws.Range("I126").Value = 552004
$ws.Range("J126").Value = 5625
$ws.Range("K126").Value = 1656012
$ws.Range("L126").Value = 16875
$ws.Range("M126").Value = -1653542
$ws.Range("N126").Value = -21815
